$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns at D:E, shifting the existing quarterly data (previously D:K) to F:M
$ws.Columns("D:E").Insert()

# Copy formatting from column F (the original column D, now shifted) onto the new D:E columns
# so the new cells keep the same number formats/fonts as the rest of the table.
$ws.Range("F7:F35").Copy() | Out-Null
$ws.Range("D7:E35").PasteSpecial(-4122) | Out-Null
$ws.Range("F38:F77").Copy() | Out-Null
$ws.Range("D38:E77").PasteSpecial(-4122) | Out-Null
$ws.Range("F80:F102").Copy() | Out-Null
$ws.Range("D80:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the two new columns (D = newest quarter, E = second-newest quarter) with the newly reported figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 353200
$ws.Cells.Item(8, 5).Value = 391700
$ws.Cells.Item(9, 4).Value = 75000
$ws.Cells.Item(9, 5).Value = 78900
$ws.Cells.Item(10, 4).Value = 278200
$ws.Cells.Item(10, 5).Value = 312800
$ws.Cells.Item(12, 4).Value = 175400
$ws.Cells.Item(12, 5).Value = 161400
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = "NA"
$ws.Cells.Item(15, 4).Value = 7600
$ws.Cells.Item(15, 5).Value = 7600
$ws.Cells.Item(17, 4).Value = 391400
$ws.Cells.Item(17, 5).Value = 407400
$ws.Cells.Item(18, 4).Value = -38200
$ws.Cells.Item(18, 5).Value = -15700
$ws.Cells.Item(20, 4).Value = 2600
$ws.Cells.Item(20, 5).Value = 8500
$ws.Cells.Item(21, 4).Value = -10800
$ws.Cells.Item(21, 5).Value = 19100
$ws.Cells.Item(22, 4).Value = 7700
$ws.Cells.Item(22, 5).Value = 12100
$ws.Cells.Item(23, 4).Value = -43300
$ws.Cells.Item(23, 5).Value = -19400
$ws.Cells.Item(24, 4).Value = -39700
$ws.Cells.Item(24, 5).Value = -6800
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = -3700
$ws.Cells.Item(26, 5).Value = -12600
$ws.Cells.Item(27, 4).Value = -3700
$ws.Cells.Item(27, 5).Value = -12600
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -2600
$ws.Cells.Item(32, 5).Value = -8500
$ws.Cells.Item(33, 4).Value = -3700
$ws.Cells.Item(33, 5).Value = -12600
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = -3700
$ws.Cells.Item(35, 5).Value = -12600
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 494000
$ws.Cells.Item(41, 5).Value = 882200
$ws.Cells.Item(42, 4).Value = 590300
$ws.Cells.Item(42, 5).Value = 561200
$ws.Cells.Item(43, 4).Value = 342600
$ws.Cells.Item(43, 5).Value = 384300
$ws.Cells.Item(44, 4).Value = 530900
$ws.Cells.Item(44, 5).Value = 508500
$ws.Cells.Item(45, 4).Value = 98400
$ws.Cells.Item(45, 5).Value = 71700
$ws.Cells.Item(46, 4).Value = 2056200
$ws.Cells.Item(46, 5).Value = 2407900
$ws.Cells.Item(47, 4).Value = 235900
$ws.Cells.Item(47, 5).Value = 204900
$ws.Cells.Item(48, 4).Value = 948700
$ws.Cells.Item(48, 5).Value = 924000
$ws.Cells.Item(49, 4).Value = 688800
$ws.Cells.Item(49, 5).Value = 691700
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 497500
$ws.Cells.Item(52, 5).Value = 466800
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 4427100
$ws.Cells.Item(54, 5).Value = 4695300
$ws.Cells.Item(57, 4).Value = 208100
$ws.Cells.Item(57, 5).Value = 180600
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 374200
$ws.Cells.Item(59, 4).Value = 315100
$ws.Cells.Item(59, 5).Value = 283800
$ws.Cells.Item(60, 4).Value = 523200
$ws.Cells.Item(60, 5).Value = 838700
$ws.Cells.Item(61, 4).Value = 830400
$ws.Cells.Item(61, 5).Value = 826100
$ws.Cells.Item(62, 4).Value = 105500
$ws.Cells.Item(62, 5).Value = 115500
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 1459200
$ws.Cells.Item(66, 5).Value = 1780300
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = -1694100
$ws.Cells.Item(72, 5).Value = -1690500
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 2967900
$ws.Cells.Item(76, 5).Value = 2914900
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = -3700
$ws.Cells.Item(81, 5).Value = -12600
$ws.Cells.Item(83, 4).Value = 24700
$ws.Cells.Item(83, 5).Value = 26300
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 49400
$ws.Cells.Item(89, 5).Value = 47200
$ws.Cells.Item(91, 4).Value = -42100
$ws.Cells.Item(91, 5).Value = -49900
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -71000
$ws.Cells.Item(94, 5).Value = 402300
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -366400
$ws.Cells.Item(100, 5).Value = 5200
$ws.Cells.Item(101, 4).Value = -300
$ws.Cells.Item(101, 5).Value = 100
$ws.Cells.Item(102, 4).Value = -388200
$ws.Cells.Item(102, 5).Value = 454800
